$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
}

Replace-Text "2024-12-27 Friday" "2024-12-28 Saturday"

Replace-Text "46×60=" "82×97="
Replace-Text "70×93=" "44×77="
Replace-Text "43×80=" "52×58="
Replace-Text "88×90=" "47×50="
Replace-Text "89×66=" "54×97="
Replace-Text "94×40=" "84×85="
Replace-Text "67×99=" "88×41="
Replace-Text "32×21=" "42×71="
Replace-Text "13×14=" "19×81="
Replace-Text "99×47=" "86×78="
Replace-Text "81×99=" "48×31="
Replace-Text "83×11=" "30×95="
Replace-Text "33×49=" "24×44="
Replace-Text "79×19=" "79×77="
Replace-Text "57×25=" "33×76="
Replace-Text "80×92=" "90×77="
Replace-Text "94×67=" "74×29="
Replace-Text "72×33=" "66×25="
Replace-Text "91×98=" "59×58="
Replace-Text "69×69=" "25×49="
Replace-Text "20×61=" "20×55="
Replace-Text "26×61=" "76×76="
Replace-Text "51×76=" "38×82="
Replace-Text "73×26=" "46×32="
Replace-Text "42×51=" "84×59="
